$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1819176666666666
$ws.Range("H2").Value = 0.5457529999999999
$ws.Range("I2").Value = 0.2278676299892611
$ws.Range("J2").Value = 0.2278676299892611
$ws.Range("M2").Value = 0.003643333333333333
$ws.Range("N2").Value = 0.01093
$ws.Range("O2").Value = 0.002177035403614994
$ws.Range("P2").Value = 0.002177035403614994
$ws.Range("Q2").Value = 0.0006627866988888889
$ws.Range("R2").Value = 0.00596508029
$ws.Range("S2").Value = 0.0004960758978244632
$ws.Range("T2").Value = 0.0004960758978244632

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1819176666666666
$ws.Range("H3").Value = 0.5457529999999999
$ws.Range("I3").Value = 0.2278676299892611
$ws.Range("J3").Value = 0.2278676299892611
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.669886333333333
$ws.Range("N3").Value = 5.009659
$ws.Range("O3").Value = 0.997822964596385
$ws.Range("P3").Value = 0.997822964596385
$ws.Range("Q3").Value = 0.3037818253585555
$ws.Range("R3").Value = 2.734036428227
$ws.Range("S3").Value = 0.2273715540914366
$ws.Range("T3").Value = 0.2273715540914367

# Row 4
$ws.Range("H4").Value = 0.8572740000000001
$ws.Range("I4").Value = 0.3579366391598651
$ws.Range("J4").Value = 0.3579366391598652
$ws.Range("M4").Value = 0.003643333333333333
$ws.Range("N4").Value = 0.01093
$ws.Range("O4").Value = 0.002177035403614994
$ws.Range("P4").Value = 0.002177035403614994
$ws.Range("Q4").Value = 0.001041111646666667
$ws.Range("R4").Value = 0.009370004820000001
$ws.Range("S4").Value = 0.0007792407357019915
$ws.Range("T4").Value = 0.0007792407357019917

# Row 5
$ws.Range("H5").Value = 0.8572740000000001
$ws.Range("I5").Value = 0.3579366391598651
$ws.Range("J5").Value = 0.3579366391598652
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.669886333333333
$ws.Range("N5").Value = 5.009659
$ws.Range("O5").Value = 0.997822964596385
$ws.Range("P5").Value = 0.997822964596385
$ws.Range("Q5").Value = 0.4771833788406667
$ws.Range("R5").Value = 4.294650409566001
$ws.Range("S5").Value = 0.3571573984241631
$ws.Range("T5").Value = 0.3571573984241632

# Row 6
$ws.Range("G6").Value = 0.3306723333333333
$ws.Range("H6").Value = 0.992017
$ws.Range("I6").Value = 0.4141957308508737
$ws.Range("J6").Value = 0.4141957308508737
$ws.Range("M6").Value = 0.003643333333333333
$ws.Range("N6").Value = 0.01093
$ws.Range("O6").Value = 0.002177035403614994
$ws.Range("P6").Value = 0.002177035403614994
$ws.Range("Q6").Value = 0.001204749534444445
$ws.Range("R6").Value = 0.01084274581
$ws.Range("S6").Value = 0.0009017187700885393
$ws.Range("T6").Value = 0.0009017187700885394

# Row 7
$ws.Range("G7").Value = 0.3306723333333333
$ws.Range("H7").Value = 0.992017
$ws.Range("I7").Value = 0.4141957308508737
$ws.Range("J7").Value = 0.4141957308508737
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.669886333333333
$ws.Range("N7").Value = 5.009659
$ws.Range("O7").Value = 0.997822964596385
$ws.Range("P7").Value = 0.997822964596385
$ws.Range("Q7").Value = 0.5521852102447778
$ws.Range("R7").Value = 4.969666892203001
$ws.Range("S7").Value = 0.4132940120807851
$ws.Range("T7").Value = 0.4132940120807852
